# Converts a "RRGGBB" hex string into the VBA-style BGR-packed long that
# PowerPoint's RGB color properties expect (R + G*256 + B*65536).
function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1) Table style swap on the 3 tables (slides 14, 15, 16) -------------
$oldStyleId = "{618C6CA9-A5C4-4A07-880D-C309BF49E5E0}"
$newStyleId = "{19738D0E-CAA7-497F-9176-7D5F7BD384D6}"

foreach ($slideIdx in 14,15,16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Swap the two embedded themes --------------------------------------
# ppt/theme/theme1.xml (used by the Notes Master) currently holds the
# "Office Theme" / "Office" colour scheme; ppt/theme/theme2.xml (used by the
# Slide Master) currently holds the "Integral" / "Red Violet" colour scheme.
# The edit swaps their contents, so theme1 ends up with the Red
# Violet/Integral colours and theme2 ends up with the Office colours.
# fontScheme/fmtScheme are identical between the two themes already, so only
# the 12 clrScheme colour slots need to move (order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink).

$officeColors = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")
$integralColors = @("000000","FFFFFF","454551","D8D9DC","E32D91","C830CC","4EA6DC","4775E7","8971E1","D54773","6B9F25","8C8C8C")

$notesTheme = $p.NotesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesTheme.Item($i).RGB = HexToRgb $integralColors[$i - 1]
}

$slideTheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $slideTheme.Item($i).RGB = HexToRgb $officeColors[$i - 1]
}
